$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) AgreementInfo: remove the "DOT" test-case row (row 4), shifting
#    the "ROW" row up from row 5 to row 4.
# ------------------------------------------------------------------
$wsAgreement = $wb.Worksheets.Item("AgreementInfo")
$wsAgreement.Rows("4:4").Delete()
[void]$wsAgreement.Range("B7").Select()

# ------------------------------------------------------------------
# 2) Add the new "PayeeInfo" worksheet as the last sheet (after
#    AgreementInfo). It becomes the active sheet.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPayee = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsPayee.Name = "PayeeInfo"

# Copy header-row formatting (bold + yellow fill) from AgreementInfo's
# row 2 (A2:C2 => styles 21,22,21) so the new sheet re-uses the same
# cell styles as the rest of the workbook.
$wsAgreement.Range("A2:B2").Copy()
$wsPayee.Range("A2:B2").PasteSpecial(-4122)
$wsAgreement.Range("A2").Copy()
$wsPayee.Range("C2").PasteSpecial(-4122)

# Copy the date-format style used on AgreementInfo's B3 cell (style 12)
# onto PayeeInfo's B3 cell.
$wsAgreement.Range("B3").Copy()
$wsPayee.Range("B3").PasteSpecial(-4122)

$wsPayee.Application.CutCopyMode = $false

# Populate values in the exact order the strings were introduced so
# that the shared-string table ends up in the same sequence:
#   AddPayeeInformationALT, LandownerName,
#   "Parcel #: test1, Grantor Name: , County PID:", AvailableTract
$wsPayee.Range("A3").Value = "AddPayeeInformationALT"
$wsPayee.Range("B2").Value = "LandownerName"
$wsPayee.Range("C3").Value = "Parcel #: test1, Grantor Name: , County PID:"
$wsPayee.Range("C2").Value = "AvailableTract"

$wsPayee.Range("A1").Value = "String"
$wsPayee.Range("B1").Value = "String"
$wsPayee.Range("C1").Value = "String"
$wsPayee.Range("A2").Value = "Title"
$wsPayee.Range("B3").Value = "Long Form Renewable"

# Column widths (best-fit approximation; the underlying engine quantizes
# widths to 1/6-character increments, so these are the closest values
# to the originally recorded best-fit widths of 28.85546875 / 20.85546875 / 40).
$wsPayee.Columns("A").ColumnWidth = 28
$wsPayee.Columns("B").ColumnWidth = 20
$wsPayee.Columns("C").ColumnWidth = 39.166666666666664

[void]$wsPayee.Range("C6").Select()
